$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 18:05"

# Update case counts for several provincias whose totals were refreshed
$ws.Range("B4").Value = 66302
$ws.Range("E4").Value = 8863

$ws.Range("B5").Value = 55921
$ws.Range("E5").Value = 5956

$ws.Range("B6").Value = 18512
$ws.Range("E6").Value = 1953

$ws.Range("B7").Value = 16648
$ws.Range("E7").Value = 2898

$ws.Range("B9").Value = 12458

$ws.Range("B14").Value = 5503

$ws.Range("B16").Value = 5178

$ws.Range("E20").Value = 351

$ws.Range("B32").Value = 2369
$ws.Range("E32").Value = 317

# Rows 33 and 34 swap identity (Soria <-> Gran Canaria) because the
# underlying shared-string order changed; the row that used to show
# "Soria" now shows "Gran Canaria" with refreshed figures, and the row
# that used to show "Gran Canaria" now shows "Soria" with the previous
# Soria figures.
$ws.Range("A33").Value = "Gran Canaria"
$ws.Range("B33").Value = 2294
$ws.Range("C33").Value = 1524
$ws.Range("D33").Value = 614
$ws.Range("E33").Value = 153

$ws.Range("A34").Value = "Soria"
$ws.Range("B34").Value = 2290
$ws.Range("C34").Value = 397
$ws.Range("D34").Value = 1774
$ws.Range("E34").Value = 119
